# Edit script: "Updated the story and lore"
# Applies the diff turning the "Today" section into "Wrecked Lands" (split into
# its own paragraph), removes the stray _GoBack bookmark from the robots
# paragraph, and appends the new "Red Desert" / "Ending" lore section at the
# end of the document (moving the _GoBack bookmark there too).

$d = $word.ActiveDocument

# --- 1. Split the "Today: <break><break>The human society..." paragraph into
#        a bold/underlined "Wrecked Lands:" heading paragraph followed by a
#        plain paragraph with the body text. ---
$headingPara = $d.Paragraphs(5)
$headingText = $headingPara.Range.Text
if ($headingText -notmatch "Today") {
    Write-Host "WARNING: paragraph 5 did not contain 'Today' as expected:" $headingText
}
$xmlToday = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:br/></w:r><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>Wrecked Lands</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">The human society has crumbled to war and the </w:t></w:r><w:r><w:t>zombies</w:t></w:r><w:r><w:t>, hungry to eat them alive. With everything destroyed, everyone is fighting to survive. Those remaining seek refuge in abandoned buildings, leftover houses and a constant fight against hunger. The food is hoarded by people and is hard to find and the danger of encountering the mutants always lingers around.</w:t></w:r></w:p>
'@
$headingPara.Range.InsertXML($xmlToday)

# --- 2. Update the "robots / Map" paragraph: give its paragraph mark the
#        bold/underline run formatting, drop the _GoBack bookmark that used to
#        sit mid-paragraph, and append the new trailing paragraphs (blank
#        heading-styled paragraph, "Red Desert:", the two desert paragraphs,
#        and the "Ending:" paragraph that now owns the _GoBack bookmark). ---
$robotsPara = $d.Paragraphs(7)
$robotsText = $robotsPara.Range.Text
if ($robotsText -notmatch "Map:") {
    Write-Host "WARNING: paragraph 7 did not contain 'Map:' as expected:" $robotsText
}
$xmlTail = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main"><w:pPr><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:t>The robots that were used in the war, run awry and shoot people at sight. With the ongoing zombie o</w:t></w:r><w:r><w:t>utbreak, people have to not only</w:t></w:r><w:r><w:t xml:space="preserve"> fight the zombies but also encounter these deadly war robots. </w:t></w:r><w:r><w:br/></w:r><w:r><w:br/></w:r><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>Map:</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0"><wp:extent cx="5935980" cy="2994660"/><wp:effectExtent l="0" t="0" r="7620" b="0"/><wp:docPr id="1" name="Picture 1"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="0" name="Picture 1"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId4" cstate="print"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="5935980" cy="2994660"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln><a:noFill/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main"><w:pPr><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main"><w:pPr><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:lastRenderedPageBreak/><w:t>Red Desert:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">The war and widespread disease and destruction has not even spared the desert. The radiation and the contamination with the &#8220;Forever Elixir&#8221; has caused the creatures to mutate. These creatures are abundant in the region and now pose a threat to the citizens who want to flee the town. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">The desert area also has a distribution center (factory) for this Elixir which is destroyed by the war and has contaminated the land. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main"><w:pPr><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>Ending:</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:br/></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
$robotsPara.Range.InsertXML($xmlTail)

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
